$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 51 cell update(s) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 228.25
$ws.Range("J19").Value = 220.125
$ws.Range("L19").Value = 220.125
$ws.Range("N19").Value = -570.125
$ws.Range("H33").Value = 484.5
$ws.Range("J33").Value = 816
$ws.Range("L33").Value = 816
$ws.Range("N33").Value = -1274
$ws.Range("H69").Value = 7681
$ws.Range("J69").Value = 9014
$ws.Range("L69").Value = 27042
$ws.Range("N69").Value = -28790
$ws.Range("H72").Value = 7681
$ws.Range("J72").Value = 9014
$ws.Range("L72").Value = 81126
$ws.Range("N72").Value = -89862
$ws.Range("H92").Value = 1082.9286
$ws.Range("I92").Value = 1114.1818
$ws.Range("J92").Value = 968.3333
$ws.Range("K92").Value = 1114.1818
$ws.Range("L92").Value = 968.3333
$ws.Range("M92").Value = 133.8181999999999
$ws.Range("N92").Value = -3464.3333
$ws.Range("H98").Value = 1150.7
$ws.Range("I98").Value = 1150.7
$ws.Range("K98").Value = 1150.7
$ws.Range("M98").Value = 347.3
$ws.Range("H113").Value = 2991.25
$ws.Range("I113").Value = 2991.25
$ws.Range("K113").Value = 2991.25
$ws.Range("M113").Value = 262.75
$ws.Range("H122").Value = 1150.7
$ws.Range("I122").Value = 1150.7
$ws.Range("K122").Value = 3452.1
$ws.Range("M122").Value = -1002.1
$ws.Range("H132").Value = 1986.409
$ws.Range("I132").Value = 1914.3334
$ws.Range("K132").Value = 5743.0002
$ws.Range("M132").Value = -3213.0002
$ws.Range("H135").Value = 1015.73914
$ws.Range("I135").Value = 498.52942
$ws.Range("K135").Value = 4486.76478
$ws.Range("M135").Value = -1951.76478
$ws.Range("H137").Value = 2012.16
$ws.Range("I137").Value = 1101.1538
$ws.Range("K137").Value = 3303.4614
$ws.Range("M137").Value = -753.4614000000001
$ws.Range("H138").Value = 3970.3606
$ws.Range("J138").Value = 4328.037
$ws.Range("L138").Value = 12984.111
$ws.Range("N138").Value = -23264.111

# --- Sheet ARM: 23 cell update(s) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16091.232
$ws.Range("I32").Value = 6888.5
$ws.Range("J32").Value = 26709.77
$ws.Range("K32").Value = 6888.5
$ws.Range("L32").Value = 26709.77
$ws.Range("M32").Value = -6601.5
$ws.Range("N32").Value = -27283.77
$ws.Range("H74").Value = 3601.2354
$ws.Range("I74").Value = 1128.125
$ws.Range("K74").Value = 1128.125
$ws.Range("M74").Value = -254.125
$ws.Range("H77").Value = 3601.2354
$ws.Range("I77").Value = 1128.125
$ws.Range("K77").Value = 5640.625
$ws.Range("M77").Value = -1272.625
$ws.Range("H132").Value = 1249.2285
$ws.Range("I132").Value = 840.0909
$ws.Range("K132").Value = 2520.2727
$ws.Range("M132").Value = 9.727300000000014
$ws.Range("H140").Value = 107597.8
$ws.Range("J140").Value = 107597.8
$ws.Range("L140").Value = 107597.8
$ws.Range("N140").Value = -117957.8

# --- Sheet BSM: 28 cell update(s) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1289.2307
$ws.Range("I20").Value = 1224
$ws.Range("J20").Value = 1648
$ws.Range("K20").Value = 1224
$ws.Range("L20").Value = 1648
$ws.Range("M20").Value = -977
$ws.Range("N20").Value = -2142
$ws.Range("H86").Value = 3565.875
$ws.Range("I86").Value = 3552.4285
$ws.Range("J86").Value = 3576.3333
$ws.Range("K86").Value = 3552.4285
$ws.Range("L86").Value = 3576.3333
$ws.Range("M86").Value = -2429.4285
$ws.Range("N86").Value = -5822.3333
$ws.Range("H89").Value = 3565.875
$ws.Range("I89").Value = 3552.4285
$ws.Range("J89").Value = 3576.3333
$ws.Range("K89").Value = 17762.1425
$ws.Range("L89").Value = 17881.6665
$ws.Range("M89").Value = -12146.1425
$ws.Range("N89").Value = -29113.6665
$ws.Range("H134").Value = 2342.9285
$ws.Range("I134").Value = 2124.16
$ws.Range("J134").Value = 4166
$ws.Range("K134").Value = 6372.48
$ws.Range("L134").Value = 12498
$ws.Range("M134").Value = -3837.48
$ws.Range("N134").Value = -17568

# --- Sheet CRP: 42 cell update(s) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 13428.571
$ws.Range("J4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("N4").Value = -3224
$ws.Range("H86").Value = 12677.143
$ws.Range("J86").Value = 15099
$ws.Range("L86").Value = 15099
$ws.Range("N86").Value = -17345
$ws.Range("H89").Value = 12677.143
$ws.Range("J89").Value = 15099
$ws.Range("L89").Value = 75495
$ws.Range("N89").Value = -86727
$ws.Range("H99").Value = 11971.071
$ws.Range("J99").Value = 13453.333
$ws.Range("L99").Value = 13453.333
$ws.Range("N99").Value = -16449.333
$ws.Range("H105").Value = 3476.182
$ws.Range("I105").Value = 1068
$ws.Range("J105").Value = 5483
$ws.Range("K105").Value = 1068
$ws.Range("L105").Value = 5483
$ws.Range("M105").Value = 679
$ws.Range("N105").Value = -8977
$ws.Range("H122").Value = 6236.76
$ws.Range("I122").Value = 6295.9287
$ws.Range("J122").Value = 6161.4546
$ws.Range("K122").Value = 18887.7861
$ws.Range("L122").Value = 18484.3638
$ws.Range("M122").Value = -16437.7861
$ws.Range("N122").Value = -23384.3638
$ws.Range("H126").Value = 11971.071
$ws.Range("J126").Value = 13453.333
$ws.Range("L126").Value = 40359.999
$ws.Range("N126").Value = -45299.999
$ws.Range("H132").Value = 5100.2
$ws.Range("I132").Value = 3240
$ws.Range("K132").Value = 9720
$ws.Range("M132").Value = -7190
$ws.Range("H141").Value = 168331.67
$ws.Range("J141").Value = 168331.67
$ws.Range("L141").Value = 168331.67
$ws.Range("N141").Value = -178691.67

# --- Sheet CUL: 16 cell update(s) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 93281.73
$ws.Range("J55").Value = 3528.4285
$ws.Range("L55").Value = 10585.2855
$ws.Range("N55").Value = -10939.2855
$ws.Range("H117").Value = 1707.5
$ws.Range("J117").Value = 2781.6667
$ws.Range("L117").Value = 8345.000100000001
$ws.Range("N117").Value = -15229.0001
$ws.Range("H139").Value = 4788.0625
$ws.Range("I139").Value = 2240.182
$ws.Range("K139").Value = 6720.545999999999
$ws.Range("M139").Value = -1580.545999999999
$ws.Range("H140").Value = 3247.4443
$ws.Range("I140").Value = 3247.4443
$ws.Range("K140").Value = 9742.332900000001
$ws.Range("M140").Value = -4562.332900000001

# --- Sheet GSM: 32 cell update(s) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7999.4
$ws.Range("I70").Value = 7997
$ws.Range("J70").Value = 8000
$ws.Range("K70").Value = 7997
$ws.Range("L70").Value = 8000
$ws.Range("M70").Value = -7727
$ws.Range("N70").Value = -8540
$ws.Range("H73").Value = 7999.4
$ws.Range("I73").Value = 7997
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 7997
$ws.Range("L73").Value = 8000
$ws.Range("M73").Value = -7061
$ws.Range("N73").Value = -9872
$ws.Range("H122").Value = 396689.66
$ws.Range("I122").Value = 70390.47
$ws.Range("J122").Value = 773188.7
$ws.Range("K122").Value = 211171.41
$ws.Range("L122").Value = 2319566.1
$ws.Range("M122").Value = -208721.41
$ws.Range("N122").Value = -2324466.1
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 2494.2856
$ws.Range("I132").Value = 2007.3636
$ws.Range("J132").Value = 4279.6665
$ws.Range("K132").Value = 6022.0908
$ws.Range("L132").Value = 12838.9995
$ws.Range("M132").Value = -3492.0908
$ws.Range("N132").Value = -17898.9995

# --- Sheet LTW: 36 cell update(s) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1774.1111
$ws.Range("I7").Value = 1567.1428
$ws.Range("J7").Value = 2498.5
$ws.Range("K7").Value = 1567.1428
$ws.Range("L7").Value = 2498.5
$ws.Range("M7").Value = -1455.1428
$ws.Range("N7").Value = -2722.5
$ws.Range("H82").Value = 1574.909
$ws.Range("I82").Value = 1523.4667
$ws.Range("J82").Value = 1685.1428
$ws.Range("K82").Value = 1523.4667
$ws.Range("L82").Value = 1685.1428
$ws.Range("M82").Value = -1162.4667
$ws.Range("N82").Value = -2407.1428
$ws.Range("H85").Value = 1574.909
$ws.Range("I85").Value = 1523.4667
$ws.Range("J85").Value = 1685.1428
$ws.Range("K85").Value = 1523.4667
$ws.Range("L85").Value = 1685.1428
$ws.Range("M85").Value = -275.4666999999999
$ws.Range("N85").Value = -4181.1428
$ws.Range("H126").Value = 1774.1111
$ws.Range("I126").Value = 1567.1428
$ws.Range("J126").Value = 2498.5
$ws.Range("K126").Value = 4701.428400000001
$ws.Range("L126").Value = 7495.5
$ws.Range("M126").Value = -2231.428400000001
$ws.Range("N126").Value = -12435.5
$ws.Range("H132").Value = 5299.407
$ws.Range("I132").Value = 4849.357
$ws.Range("K132").Value = 14548.071
$ws.Range("M132").Value = -12018.071
$ws.Range("H141").Value = 70709
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# --- Sheet WVR: 8 cell update(s) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2500375
$ws.Range("J2").Value = 2500375
$ws.Range("L2").Value = 2500375
$ws.Range("N2").Value = -2500599
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

Write-Host "Applied scheduled-runner price updates across all sheets."
